# RMI files update 3.4.3
# - "About" sheet: add region label "Minnesota" in B1, bump the date in C1
# - "BIEfIE" sheet: flip the boolean control-lever value in B2 from 1 to 0

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$about.Range("B1").Value = "Minnesota"
$about.Range("C1").Value = 44840

$bie = $wb.Worksheets.Item("BIEfIE")
$bie.Range("B2").Value = 0
